$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in rows 2-5 (column A: cluster ids, column B: counts)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 57

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 40

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 34

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 33

# Remove row 6 entirely (shrinks used range/dimension to A1:B5)
$ws.Rows(6).Delete()
